$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.854.37"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.634.61"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'215.12"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "'0.0643"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "'20.00"
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("D11").Value = "'0.0782"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "'4.26"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "1.637.59"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "1.858.07"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "'0.561"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "0.0₃0766"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "'63.10"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "25.853.41"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "'194.41"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "'4.39"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").Value = "'9.93"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").Value = "'6.18"
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'1.76"
$ws.Range("E25").Value = "  -4.69%  "
$ws.Range("D26").Value = "'138.60"
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("E27").Value = "  -4.21%  "
$ws.Range("D28").Value = "'6.82"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("D29").Value = "'15.56"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "'0.0494"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'3.25"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("D34").Value = "'1.57"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("D36").Value = "'0.901"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("D38").Value = "1.125.93"
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("D39").Value = "'0.548"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "'0.0155"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "'5.51"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("D43").Value = "'99.55"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").Value = "'0.800"
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("D45").Value = "0.0₆0109"
$ws.Range("E45").Value = "  -4.82%  "
$ws.Range("D46").Value = "'55.42"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("E47").Value = "  -4.96%  "
$ws.Range("D48").Value = "'0.0503"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").Value = "'7.61"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("E51").Value = "  -0.51%  "
